# TODO.xlsx - replace the finished/old task list with the current one and
# drop the "Good" (green) highlight style that used to mark "Ignas" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all the old task rows (rows 2-17), keeping only the header row.
$ws.Rows("5:17").Delete()

# The 3 remaining data rows used the conditional "Good" cell style (green
# fill/font) on several cells - strip that formatting so the cells fall back
# to the default/normal style, matching the new plain task list look.
$ws.Range("A2:B4").ClearFormats()

# Write the new, current task list.
$ws.Range("A2").Value = "Autoscout24.de"
$ws.Range("B2").Value = "Artas"

$ws.Range("A3").Value = "Logas parserio"
$ws.Range("B3").Value = "Ignas"

$ws.Range("A4").Value = "Edit disable settings"
$ws.Range("B4").Value = "Ignas"

# The "Good" named cell style is no longer used anywhere - delete it so the
# workbook goes back to only having the standard "Normal" style.
$wb.Styles.Item("Good").Delete()

# Move the active selection to right after the new table, like a user would
# leave it after typing in the last row.
[void]$ws.Range("A5").Select()
